$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.300.34'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.279.27'
$ws.Range("E3").Value = '  -2.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.85'
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.63'
$ws.Range("E6").Value = '  -7.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.578'
$ws.Range("E8").Value = '  -3.50%  '
$ws.Range("D9").Value = '3.274.56'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  -5.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.573'
$ws.Range("E11").Value = '  -2.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.10'
$ws.Range("E12").Value = '  -5.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '666.83'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '3.797.83'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.29'
$ws.Range("E16").Value = '  -4.12%  '
$ws.Range("D17").Value = '67.158.04'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '3.258.33'
$ws.Range("E19").Value = '  -3.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.30'
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.77'
$ws.Range("E21").Value = '  -4.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.878'
$ws.Range("E22").Value = '  -3.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.33'
$ws.Range("E23").Value = '  +4.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.98'
$ws.Range("E24").Value = '  -6.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.55'
$ws.Range("E25").Value = '  -3.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.86'
$ws.Range("E26").Value = '  -3.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.65'
$ws.Range("E27").Value = '  -6.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.21'
$ws.Range("E28").Value = '  -6.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.48'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.32'
$ws.Range("E30").Value = '  -5.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.87'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '567.78'
$ws.Range("E32").Value = '  -8.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.89'
$ws.Range("E33").Value = '  -3.06%  '
$ws.Range("D34").Value = '3.758.45'
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.103'
$ws.Range("E35").Value = '  -4.17%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.40'
$ws.Range("E37").Value = '  -13.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.64'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.130'
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.63'
$ws.Range("E40").Value = '  -7.32%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.42'
$ws.Range("E41").Value = '  -4.32%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.04'
$ws.Range("E42").Value = '  -7.47%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0661'
$ws.Range("E43").Value = '  -6.99%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.326'
$ws.Range("E44").Value = '  -5.90%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.21'
$ws.Range("E45").Value = '  -5.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0400'
$ws.Range("E46").Value = '  -5.82%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.59'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("E48").Value = '  -2.95%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.74'
$ws.Range("E51").Value = '  -4.76%  '
